$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the USD amount recorded in T2 for this deposit row
$ws.Range("T2").Value = 205165

# Move the active selection to T3 (matches the saved cursor position)
$ws.Range("T3").Select()
